# Add a new "About" sheet at the front of the workbook, describing the
# framework the file is drawn from (name + description), matching the
# "Add about sheet from which name is drawn" commit.

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current first sheet so it becomes tab 1.
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "About"

# Header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Description"

# Data row.
$ws.Range("A2").Value = "SIR"
$ws.Range("B2").Value = "The SIR model"

# Formatting: bold header row, top-aligned data row.
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A2:B2").VerticalAlignment = -4160

# Leave the selection on B2, matching the authored selection state.
[void]$ws.Range("B2").Select()
